$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.434825
$ws.Range("H2").Value = 4.304475
$ws.Range("I2").Value = 0.3232227003447979
$ws.Range("J2").Value = 0.3232227003447979
$ws.Range("M2").Value = 1.699817666666666
$ws.Range("N2").Value = 5.099453
$ws.Range("O2").Value = 0.748079722752454
$ws.Range("P2").Value = 0.748079722752454
$ws.Range("Q2").Value = 2.438940883574999
$ws.Range("R2").Value = 21.950467952175
$ws.Range("S2").Value = 0.2417963480612359
$ws.Range("T2").Value = 0.241796348061236
# Row 3
$ws.Range("G3").Value = 1.434825
$ws.Range("H3").Value = 4.304475
$ws.Range("I3").Value = 0.3232227003447979
$ws.Range("J3").Value = 0.3232227003447979
$ws.Range("M3").Value = 0.5724236666666667
$ws.Range("N3").Value = 1.717271
$ws.Range("O3").Value = 0.2519202772475459
$ws.Range("P3").Value = 0.2519202772475459
$ws.Range("Q3").Value = 0.821327787525
$ws.Range("R3").Value = 7.391950087725
$ws.Range("S3").Value = 0.08142635228356193
$ws.Range("T3").Value = 0.08142635228356195
# Row 4
$ws.Range("I4").Value = 0.2419118765677838
$ws.Range("J4").Value = 0.2419118765677838
$ws.Range("M4").Value = 1.699817666666666
$ws.Range("N4").Value = 5.099453
$ws.Range("O4").Value = 0.748079722752454
$ws.Range("P4").Value = 0.748079722752454
$ws.Range("Q4").Value = 1.825393963215222
$ws.Range("R4").Value = 16.428545668937
$ws.Range("S4").Value = 0.1809693695533535
$ws.Range("T4").Value = 0.1809693695533535
# Row 5
$ws.Range("I5").Value = 0.2419118765677838
$ws.Range("J5").Value = 0.2419118765677838
$ws.Range("M5").Value = 0.5724236666666667
$ws.Range("N5").Value = 1.717271
$ws.Range("O5").Value = 0.2519202772475459
$ws.Range("P5").Value = 0.2519202772475459
$ws.Range("Q5").Value = 0.6147122282732222
$ws.Range("R5").Value = 5.532410054459
$ws.Range("S5").Value = 0.06094250701443019
$ws.Range("T5").Value = 0.06094250701443019
# Row 6
$ws.Range("G6").Value = 0.466371
$ws.Range("H6").Value = 1.399113
$ws.Range("I6").Value = 0.1050592887512441
$ws.Range("J6").Value = 0.1050592887512441
$ws.Range("M6").Value = 1.699817666666666
$ws.Range("N6").Value = 5.099453
$ws.Range("O6").Value = 0.748079722752454
$ws.Range("P6").Value = 0.748079722752454
$ws.Range("Q6").Value = 0.792745665021
$ws.Range("R6").Value = 7.134710985189
$ws.Range("S6").Value = 0.07859272360160066
$ws.Range("T6").Value = 0.07859272360160066
# Row 7
$ws.Range("G7").Value = 0.466371
$ws.Range("H7").Value = 1.399113
$ws.Range("I7").Value = 0.1050592887512441
$ws.Range("J7").Value = 0.1050592887512441
$ws.Range("M7").Value = 0.5724236666666667
$ws.Range("N7").Value = 1.717271
$ws.Range("O7").Value = 0.2519202772475459
$ws.Range("P7").Value = 0.2519202772475459
$ws.Range("Q7").Value = 0.266961797847
$ws.Range("R7").Value = 2.402656180623
$ws.Range("S7").Value = 0.02646656514964338
$ws.Range("T7").Value = 0.02646656514964338
# Row 8
$ws.Range("G8").Value = 0.4851976666666666
$ws.Range("H8").Value = 1.455593
$ws.Range("I8").Value = 0.1093003676552856
$ws.Range("J8").Value = 0.1093003676552856
$ws.Range("M8").Value = 1.699817666666666
$ws.Range("N8").Value = 5.099453
$ws.Range("O8").Value = 0.748079722752454
$ws.Range("P8").Value = 0.748079722752454
$ws.Range("Q8").Value = 0.8247475656254443
$ws.Range("R8").Value = 7.422728090628999
$ws.Range("S8").Value = 0.08176538873230732
$ws.Range("T8").Value = 0.08176538873230733
# Row 9
$ws.Range("G9").Value = 0.4851976666666666
$ws.Range("H9").Value = 1.455593
$ws.Range("I9").Value = 0.1093003676552856
$ws.Range("J9").Value = 0.1093003676552856
$ws.Range("M9").Value = 0.5724236666666667
$ws.Range("N9").Value = 1.717271
$ws.Range("O9").Value = 0.2519202772475459
$ws.Range("P9").Value = 0.2519202772475459
$ws.Range("Q9").Value = 0.2777386274114444
$ws.Range("R9").Value = 2.499647646703
$ws.Range("S9").Value = 0.02753497892297824
$ws.Range("T9").Value = 0.02753497892297824
# Row 10
$ws.Range("G10").Value = 0.9788520000000002
$ws.Range("H10").Value = 2.936556
$ws.Range("I10").Value = 0.2205057666808887
$ws.Range("J10").Value = 0.2205057666808887
$ws.Range("M10").Value = 1.699817666666666
$ws.Range("N10").Value = 5.099453
$ws.Range("O10").Value = 0.748079722752454
$ws.Range("P10").Value = 0.748079722752454
$ws.Range("Q10").Value = 1.663869922652
$ws.Range("R10").Value = 14.974829303868
$ws.Range("S10").Value = 0.1649558928039566
$ws.Range("T10").Value = 0.1649558928039566
# Row 11
$ws.Range("G11").Value = 0.9788520000000002
$ws.Range("H11").Value = 2.936556
$ws.Range("I11").Value = 0.2205057666808887
$ws.Range("J11").Value = 0.2205057666808887
$ws.Range("M11").Value = 0.5724236666666667
$ws.Range("N11").Value = 1.717271
$ws.Range("O11").Value = 0.2519202772475459
$ws.Range("P11").Value = 0.2519202772475459
$ws.Range("Q11").Value = 0.5603180509640001
$ws.Range("R11").Value = 5.042862458676001
$ws.Range("S11").Value = 0.05554987387693216
$ws.Range("T11").Value = 0.05554987387693216
